$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.589.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.08%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.881.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.77%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.026'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +3.81%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.52%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.025'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.70%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5161'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.35%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3983'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.59%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08406'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.05%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.118'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.37%  '

# Row 11
$ws.Range("B11").Value = 'OKB'
$ws.Range("C11").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.25'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.41%  '

# Row 12
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.295'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.56%  '

# Row 13
$ws.Range("B13").Value = 'Solana'
$ws.Range("C13").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.63'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.61%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.827.96'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.11%  '

# Row 15
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.269'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.75%  '

# Row 16
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.026'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.40%  '

# Row 17
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001115'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.70%  '

# Row 18
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.43'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.91%  '

# Row 19
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06789'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.95%  '

# Row 20
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.88%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.024'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.66%  '

# Row 22
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.997'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.39%  '

# Row 23
$ws.Range("B23").Value = 'WrappedBTC'
$ws.Range("C23").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.583.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.87%  '

# Row 24
$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.21'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.16%  '

# Row 25
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.299'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.47%  '

# Row 26
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.93'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.88%  '

# Row 27
$ws.Range("B27").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C27").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.042.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.29%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.57%  '

# Row 29
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.385'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.91%  '

# Row 30
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '128.05'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.81%  '

# Row 31
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1055'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.05%  '

# Row 32
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.044'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.21%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.846'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.18%  '

# Row 34
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.654'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.00%  '

# Row 35
$ws.Range("B35").Value = 'VeChain'
$ws.Range("C35").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02445'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.48%  '

# Row 36
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06529'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.99%  '

# Row 37
$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2200'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.11%  '

# Row 38
$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.986'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.55%  '

# Row 39
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.281'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +6.18%  '

# Row 40
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.194'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.53%  '

# Row 41
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6471'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.28%  '

# Row 42
$ws.Range("B42").Value = 'InternetComputer(DFINITY)'
$ws.Range("C42").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.064'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.21%  '

# Row 43
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.31'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.01%  '

# Row 44
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6065'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.65%  '

# Row 45
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.07'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.29%  '

# Row 46
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.749'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.79%  '

# Row 47
$ws.Range("B47").Value = 'WEMIXTOKEN'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.230'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.19%  '

# Row 48
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.007'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.34%  '

# Row 49
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.215'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.26%  '

# Row 50
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.11%  '

# Row 51
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06869'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.59%  '
